$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:E) to (B:F)
$ws.Columns.Item(1).Insert()

# Copy the header formatting (bold, centered, bordered) from B1 into the new A1 cell
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Cells.Item(1, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set the new header cell text
$ws.Cells.Item(1, 1).Value = "ID"

# Row identifiers for rows 2..25 (new column A)
$ids = @(
    "Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95",
    "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22",
    "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16"
)

$r = 2
foreach ($id in $ids) {
    $ws.Cells.Item($r, 1).Value = $id
    $r = $r + 1
}
